$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.395.74"
$ws.Range("E2").Value = "  -0.05%  "
$ws.Range("D3").Value = "3.757.38"
$ws.Range("E3").Value = "  -0.65%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'595.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("D6").Value = "'169.99"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.04%  "
$ws.Range("D7").Value = "3.758.35"
$ws.Range("E7").Value = "  -0.74%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "'0.525"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.38%  "
$ws.Range("D10").Value = "'0.166"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.19%  "
$ws.Range("D11").Value = "'6.50"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.15%  "
$ws.Range("D12").Value = "'0.455"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.35%  "
$ws.Range("D13").Value = "'0.0000278"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +8.22%  "
$ws.Range("D14").Value = "'36.69"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.22%  "
$ws.Range("D15").Value = "4.385.76"
$ws.Range("E15").Value = "  -0.95%  "
$ws.Range("D16").Value = "3.766.29"
$ws.Range("E16").Value = "  -1.24%  "
$ws.Range("D17").Value = "'18.69"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.66%  "
$ws.Range("D18").Value = "67.433.35"
$ws.Range("E18").Value = "  -0.12%  "
$ws.Range("D19").Value = "'7.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.42%  "
$ws.Range("D20").Value = "'0.112"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.35%  "
$ws.Range("D21").Value = "'10.51"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.30%  "
$ws.Range("D22").Value = "'469.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.26%  "
$ws.Range("D23").Value = "'0.721"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.53%  "
$ws.Range("D24").Value = "'83.88"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.78%  "
$ws.Range("D25").Value = "'0.0000147"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -7.78%  "
$ws.Range("D26").Value = "'2.22"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.92%  "
$ws.Range("D27").Value = "'12.17"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.46%  "
$ws.Range("D28").Value = "'10.45"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.92%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.17%  "
$ws.Range("D30").Value = "'2.90"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.22%  "
$ws.Range("D31").Value = "3.905.16"
$ws.Range("E31").Value = "  -0.94%  "
$ws.Range("D32").Value = "'7.68"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.66%  "
$ws.Range("D33").Value = "'30.59"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.54%  "
$ws.Range("D34").Value = "'2.24"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.76%  "
$ws.Range("D35").Value = "'9.14"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.88%  "
$ws.Range("D36").Value = "3.719.69"
$ws.Range("E36").Value = "  -0.98%  "
$ws.Range("D37").Value = "'3.81"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.37%  "
$ws.Range("D38").Value = "'0.105"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.43%  "
$ws.Range("D39").Value = "'0.138"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.29%  "
$ws.Range("D40").Value = "'5.87"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.22%  "
$ws.Range("D41").Value = "'1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.89%  "
$ws.Range("D42").Value = "'0.999"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("D43").Value = "'0.312"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.38%  "
$ws.Range("D44").Value = "'1.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("D45").Value = "'8.74"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.09%  "
$ws.Range("D46").Value = "'1.95"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.52%  "
$ws.Range("D47").Value = "'45.86"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.73%  "
$ws.Range("D48").Value = "'400.20"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.36%  "
$ws.Range("D49").Value = "'0.000270"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.93%  "
$ws.Range("D50").Value = "'140.33"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.95%  "
$ws.Range("D51").Value = "'39.50"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.06%  "
